$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: copy formatting from row 24 (same style pattern: 7,7,7,7,3,3,5) ---
$ws.Range("A24:G24").Copy() | Out-Null
$ws.Range("A27:G27").PasteSpecial(-4122) | Out-Null

# --- Row 28: copy formatting from row 25 (same style pattern: 7,7,7,7,3,6,5) ---
$ws.Range("A25:G25").Copy() | Out-Null
$ws.Range("A28:G28").PasteSpecial(-4122) | Out-Null

# --- Row 29: copy formatting from row 26 (same style pattern: 7,7,7,7,3,6,5) ---
$ws.Range("A26:G26").Copy() | Out-Null
$ws.Range("A29:G29").PasteSpecial(-4122) | Out-Null

# --- Row 30: copy formatting from the banner row 20 (merged, style 9 across A:I) ---
$ws.Range("A20:I20").Copy() | Out-Null
$ws.Range("A30:I30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Fill in the new question rows ---
$ws.Range("A27").Value = 96
$ws.Range("B27").Value = "Unique Binary Search Tree"
$ws.Range("C27").Value = "Tree"
$ws.Range("D27").Value = "Aton"
$ws.Range("E27").Value = "On-going"
$ws.Range("F27").Value = "Medium"
$ws.Range("G27").Value = "Python"

$ws.Range("A28").Value = 108
$ws.Range("B28").Value = "Convert Sorted Array to Binary Search Tree"
$ws.Range("C28").Value = "Tree"
$ws.Range("D28").Value = "Aton"
$ws.Range("E28").Value = "On-going"
$ws.Range("F28").Value = "Easy"
$ws.Range("G28").Value = "Python"

$ws.Range("A29").Value = 669
$ws.Range("B29").Value = "Trim a Binary Search Tree"
$ws.Range("C29").Value = "Tree"
$ws.Range("D29").Value = "Aton"
$ws.Range("E29").Value = "On-going"
$ws.Range("F29").Value = "Easy"
$ws.Range("G29").Value = "Python"

# --- Row 30: new date banner, merged across A:I like the other banner rows ---
$ws.Range("A30").Value = "April 4, 2018"
$ws.Range("A30:I30").Merge() | Out-Null

# --- Update the active selection to reflect where editing left off ---
$ws.Range("H28").Select() | Out-Null
